# Rigidbody test scene scaffolding: add a new diary entry row (row 27)
# for "10 marras" documenting the start of the rigid-body chapter demo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row values
$ws.Range("A27").Value = "10 marras"

# "Kello" (time) column stores a time-of-day fraction, formatted like the
# other h:mm cells in the column (e.g. B8, B9 ... use the same number format).
$ws.Range("B27").Value = 0.75
$ws.Range("B27").NumberFormat = "h:mm"

# "Oppimisen sisältö" column, wrapped like the rest of column C.
$ws.Range("C27").Value = "Kovien kappaleiden demon aloitus"
$ws.Range("C27").WrapText = $true

# Give the new row the same visual height the other wrapped rows get
# (two wrapped lines).
$ws.Rows.Item(27).RowHeight = 29

# Scroll the view down so the new row is visible (best effort - mirrors the
# workbook being scrolled from topLeftCell A24 to A22 in the source file).
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
